$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model for threading: A holds the 0-based position index, B holds the
# direct-argument value (as text, matching the source export format which
# always carries a trailing newline).
$values = @(668, 671, 676, 682, 688, 692, 697, 715, 732, 752)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1

    $ws.Cells.Item($row, 1).Value = $i

    $cellB = $ws.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = "$($values[$i])`n"
    $cellB.Style = "Normal"
}
